# Updation in custom order for both QA and STG
# Updates the sample row on the "Sec invoice Master" sheet with a new
# FC Order ID / tracking number / carrier / invoice amount / secondary
# invoice number, matching a freshly-pulled example order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sec invoice Master")

# FC Order ID (A2) and SECONDARY INV # (I2) are derived from the order id,
# and New Invoice Amount (F2) is a currency-looking string -- all of these
# must stay as literal text, not be reinterpreted as numbers, so they are
# entered with a leading apostrophe (the same trick Excel's UI uses to force
# text entry for numeric-looking values).
$ws.Range("A2").Value = "'59071590"
$ws.Range("C2").Value = "999U684759"
$ws.Range("D2").Value = "FragilePAK"
$ws.Range("F2").Value = "'457.13"
$ws.Range("I2").Value = "'59071590+1"
